$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("P2").Value = 3.3
$ws.Range("Q2").Value = 1.34
$ws.Range("R2").Value = 1.93
$ws.Range("S2").Value = 1.87
$ws.Range("T2").Value = 1.84
$ws.Range("U2").Value = 1.96
$ws.Range("Y2").Value = 15
$ws.Range("Z2").Value = 11
$ws.Range("AA2").Value = 11
$ws.Range("AC2").Value = 20
$ws.Range("AD2").Value = 13
$ws.Range("AE2").Value = 14
$ws.Range("AF2").Value = 150
$ws.Range("AG2").Value = 48
$ws.Range("AH2").Value = 38
$ws.Range("AI2").Value = 32
$ws.Range("AK2").Value = 190
$ws.Range("AL2").Value = 130
$ws.Range("AM2").Value = 150
$ws.Range("AN2").Value = 180
$ws.Range("AO2").Value = 3.3

# Row 3
$ws.Range("Q3").Value = 2.24
$ws.Range("T3").Value = 1.61
$ws.Range("U3").Value = 1.63

# Row 4
$ws.Range("L4").Value = 1.25
$ws.Range("Q4").Value = 1.49

# Row 6
$ws.Range("F6").Value = 1.84
$ws.Range("G6").Value = 1.99
$ws.Range("L6").Value = 1.23
$ws.Range("W6").Value = 2
$ws.Range("AB6").Value = 970
$ws.Range("AF6").Value = 970
$ws.Range("AG6").Value = 970
$ws.Range("AH6").Value = 970
$ws.Range("AK6").Value = 970

# Row 7
$ws.Range("AO7").Value = 18.5

# Row 8
$ws.Range("N8").Value = 1.03

# Row 9
$ws.Range("L9").Value = 1.4
$ws.Range("Y9").Value = 10
$ws.Range("AB9").Value = 13.5
$ws.Range("AC9").Value = 8.6

# Row 10
$ws.Range("F10").Value = 9.4
$ws.Range("G10").Value = 11
$ws.Range("W10").Value = 1.1
$ws.Range("Z10").Value = 970
$ws.Range("AA10").Value = 970
$ws.Range("AB10").Value = 29
$ws.Range("AC10").Value = 970
$ws.Range("AE10").Value = 970
$ws.Range("AG10").Value = 38

# Row 11
$ws.Range("G11").Value = 5.4
$ws.Range("H11").Value = 1.68
$ws.Range("I11").Value = 1.89
$ws.Range("Q11").Value = 1.42
$ws.Range("S11").Value = 1.91
$ws.Range("V11").Value = 2.12

# Row 12
$ws.Range("F12").Value = 1.74
$ws.Range("H12").Value = 2.94
$ws.Range("Q12").Value = 1.73
$ws.Range("V12").Value = 1.26

# Row 15
$ws.Range("N15").Value = 6.8
$ws.Range("O15").Value = 1.11
$ws.Range("P15").Value = 3.45
$ws.Range("R15").Value = 1.99
$ws.Range("T15").Value = 1.4
$ws.Range("U15").Value = 3
$ws.Range("X15").Value = 55
$ws.Range("Y15").Value = 22
$ws.Range("Z15").Value = 21
$ws.Range("AA15").Value = 30
$ws.Range("AB15").Value = 38
$ws.Range("AC15").Value = 16
$ws.Range("AD15").Value = 15
$ws.Range("AE15").Value = 21
$ws.Range("AF15").Value = 42
$ws.Range("AG15").Value = 22
$ws.Range("AH15").Value = 18
$ws.Range("AI15").Value = 26
$ws.Range("AJ15").Value = 85
$ws.Range("AK15").Value = 42
$ws.Range("AL15").Value = 34
$ws.Range("AM15").Value = 50
$ws.Range("AN15").Value = 21
$ws.Range("AO15").Value = 7

# Row 16
$ws.Range("H16").Value = 2.5

# Row 17
$ws.Range("G17").Value = 7
$ws.Range("Q17").Value = 2.3
$ws.Range("W17").Value = 1.16

# Row 18
$ws.Range("N18").Value = 1.45
$ws.Range("P18").Value = 1.45
$ws.Range("S18").Value = 2.24

# Row 19
$ws.Range("H19").Value = 1.99
$ws.Range("I19").Value = 2
$ws.Range("L19").Value = 1.5
$ws.Range("M19").Value = 1.1
$ws.Range("N19").Value = 3.25
$ws.Range("V19").Value = 2
$ws.Range("AJ19").Value = 110
$ws.Range("AL19").Value = 80
$ws.Range("AM19").Value = 150

# Row 20
$ws.Range("O20").Value = 1.45
$ws.Range("Q20").Value = 2.36
$ws.Range("AE20").Value = 95
$ws.Range("AI20").Value = 100

# Row 21
$ws.Range("G21").Value = 2.64
$ws.Range("H21").Value = 3.2
$ws.Range("I21").Value = 3.25
$ws.Range("Q21").Value = 2.36
$ws.Range("S21").Value = 4.6
$ws.Range("V21").Value = 1.44
$ws.Range("W21").Value = 1.61

# Row 23
$ws.Range("H23").Value = 1.28
$ws.Range("Z23").Value = 7.2
$ws.Range("AF23").Value = 170

# Row 24
$ws.Range("L24").Value = 1.01
$ws.Range("M24").Value = 1.01
$ws.Range("N24").Value = 1.74
$ws.Range("O24").Value = 1.39
$ws.Range("P24").Value = 1.74
$ws.Range("Q24").Value = 2.18
$ws.Range("R24").Value = 1.21
$ws.Range("S24").Value = 3.4
$ws.Range("T24").Value = 1.01
$ws.Range("U24").Value = 1.01
$ws.Range("V24").Value = 1.24
$ws.Range("W24").Value = 1.98
$ws.Range("X24").Value = 1000
$ws.Range("Y24").Value = 1000
$ws.Range("Z24").Value = 1000
$ws.Range("AA24").Value = 1000
$ws.Range("AB24").Value = 1000
$ws.Range("AC24").Value = 1000
$ws.Range("AD24").Value = 1000
$ws.Range("AE24").Value = 1000
$ws.Range("AF24").Value = 1000
$ws.Range("AG24").Value = 1000
$ws.Range("AH24").Value = 1000
$ws.Range("AI24").Value = 1000
$ws.Range("AJ24").Value = 1000
$ws.Range("AK24").Value = 1000
$ws.Range("AL24").Value = 1000
$ws.Range("AM24").Value = 1000
$ws.Range("AN24").Value = 1000
$ws.Range("AO24").Value = 1000

# Row 26
$ws.Range("K26").Value = 3.2
$ws.Range("L26").Value = 1.01
$ws.Range("M26").Value = 1.13
$ws.Range("N26").Value = 2.4
$ws.Range("O26").Value = 1.61
$ws.Range("R26").Value = 1.15
$ws.Range("S26").Value = 6.4
$ws.Range("T26").Value = 2.18
$ws.Range("U26").Value = 1.67
$ws.Range("V26").Value = 1.64
$ws.Range("W26").Value = 1.36
$ws.Range("X26").Value = 7.8
$ws.Range("Y26").Value = 7.4
$ws.Range("Z26").Value = 13.5
$ws.Range("AA26").Value = 40
$ws.Range("AB26").Value = 9.800000000000001
$ws.Range("AC26").Value = 7.6
$ws.Range("AD26").Value = 13
$ws.Range("AE26").Value = 40
$ws.Range("AF26").Value = 24
$ws.Range("AG26").Value = 17.5
$ws.Range("AH26").Value = 28
$ws.Range("AI26").Value = 80
$ws.Range("AJ26").Value = 90
$ws.Range("AK26").Value = 70
$ws.Range("AL26").Value = 120
$ws.Range("AM26").Value = 240
$ws.Range("AN26").Value = 100
$ws.Range("AO26").Value = 46

# Row 28
$ws.Range("L28").Value = 1.01
$ws.Range("M28").Value = 1.01
$ws.Range("N28").Value = 1.02
$ws.Range("O28").Value = 1.3
$ws.Range("R28").Value = 1.18
$ws.Range("S28").Value = 1.01
$ws.Range("T28").Value = 1.01
$ws.Range("U28").Value = 1.01
$ws.Range("V28").Value = 1.01
$ws.Range("W28").Value = 1.01
$ws.Range("X28").Value = 1000
$ws.Range("Y28").Value = 1000
$ws.Range("Z28").Value = 1000
$ws.Range("AA28").Value = 1000
$ws.Range("AB28").Value = 1000
$ws.Range("AC28").Value = 1000
$ws.Range("AD28").Value = 1000
$ws.Range("AE28").Value = 1000
$ws.Range("AF28").Value = 1000
$ws.Range("AG28").Value = 1000
$ws.Range("AH28").Value = 1000
$ws.Range("AI28").Value = 1000
$ws.Range("AJ28").Value = 1000
$ws.Range("AK28").Value = 1000
$ws.Range("AL28").Value = 1000
$ws.Range("AM28").Value = 1000
$ws.Range("AN28").Value = 1000
$ws.Range("AO28").Value = 1000

# Row 29
$ws.Range("F29").Value = 2.48
$ws.Range("G29").Value = 4.1
$ws.Range("H29").Value = 2.38
$ws.Range("I29").Value = 2.72
$ws.Range("L29").Value = 1.01
$ws.Range("M29").Value = 1.01
$ws.Range("N29").Value = 1.53
$ws.Range("O29").Value = 1.02
$ws.Range("Q29").Value = 2.28
$ws.Range("R29").Value = 1.14
$ws.Range("S29").Value = 4.3
$ws.Range("T29").Value = 1.01
$ws.Range("U29").Value = 1.01
$ws.Range("V29").Value = 1.58
$ws.Range("W29").Value = 1.33
$ws.Range("X29").Value = 1000
$ws.Range("Y29").Value = 1000
$ws.Range("Z29").Value = 1000
$ws.Range("AA29").Value = 1000
$ws.Range("AB29").Value = 1000
$ws.Range("AC29").Value = 1000
$ws.Range("AD29").Value = 1000
$ws.Range("AE29").Value = 1000
$ws.Range("AF29").Value = 1000
$ws.Range("AG29").Value = 1000
$ws.Range("AH29").Value = 1000
$ws.Range("AI29").Value = 1000
$ws.Range("AJ29").Value = 1000
$ws.Range("AK29").Value = 1000
$ws.Range("AL29").Value = 1000
$ws.Range("AM29").Value = 1000
$ws.Range("AN29").Value = 1000
$ws.Range("AO29").Value = 1000

# Row 30
$ws.Range("G30").Value = 2.1
$ws.Range("J30").Value = 3.6
$ws.Range("L30").Value = 1.01
$ws.Range("M30").Value = 1.01
$ws.Range("N30").Value = 1.97
$ws.Range("O30").Value = 1.28
$ws.Range("R30").Value = 1.29
$ws.Range("S30").Value = 2.94
$ws.Range("T30").Value = 1.01
$ws.Range("U30").Value = 1.01
$ws.Range("V30").Value = 1.27
$ws.Range("W30").Value = 1.91
$ws.Range("X30").Value = 1000
$ws.Range("Y30").Value = 1000
$ws.Range("Z30").Value = 1000
$ws.Range("AA30").Value = 1000
$ws.Range("AB30").Value = 1000
$ws.Range("AC30").Value = 1000
$ws.Range("AD30").Value = 1000
$ws.Range("AE30").Value = 1000
$ws.Range("AF30").Value = 1000
$ws.Range("AG30").Value = 1000
$ws.Range("AH30").Value = 1000
$ws.Range("AI30").Value = 1000
$ws.Range("AJ30").Value = 1000
$ws.Range("AK30").Value = 1000
$ws.Range("AL30").Value = 1000
$ws.Range("AM30").Value = 1000
$ws.Range("AN30").Value = 1000
$ws.Range("AO30").Value = 1000
